$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Row 8 updates
$ws.Range("F8").Value = "✅"
$ws.Range("K8").Value = "✅"

# L8 looks like a date, so force text formatting before assigning,
# then reset the style so no residual number-format sticks to the cell.
$cellL8 = $ws.Range("L8")
$cellL8.NumberFormat = "@"
$cellL8.Value = "11/12/2025"
$cellL8.Style = "Normal"

$ws.Range("O8").Value = "تيم الكرد"

# Row 29 updates: B29:K29 were empty, now filled with literal text "nan"
$ws.Range("B29:K29").Value = "nan"
